$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add new header cells P1 (14) and Q1 (15), matching O1 style ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Data block B2:Q25 (24 rows x 16 cols: B..Q) ---
$data = New-Object 'double[,]' 24,16
$data[0,0] = 24.49266632314295
$data[0,1] = 18.96828366782543
$data[0,2] = 4.333588960490924
$data[0,3] = 29.21382029618269
$data[0,4] = 18.59124872237463
$data[0,5] = 2.071584505415577
$data[0,6] = 3.159216038339132
$data[0,7] = 3.138252493947185
$data[0,8] = 0
$data[0,9] = 0
$data[0,10] = 0
$data[0,11] = 0
$data[0,12] = 0
$data[0,13] = 0
$data[0,14] = 0
$data[0,15] = 14.5687144452271
$data[1,0] = 22.88573587628169
$data[1,1] = 17.85725570160992
$data[1,2] = 4.196859081625459
$data[1,3] = 27.32102133038484
$data[1,4] = 18.00587086705919
$data[1,5] = 2.076468730828436
$data[1,6] = 2.90028629522092
$data[1,7] = 2.942540537153831
$data[1,8] = 0
$data[1,9] = 0
$data[1,10] = 0
$data[1,11] = 0
$data[1,12] = 0
$data[1,13] = 0
$data[1,14] = 0
$data[1,15] = 14.3407332244305
$data[2,0] = 21.83970279736576
$data[2,1] = 17.14498019991685
$data[2,2] = 4.110482788901874
$data[2,3] = 26.09414241005118
$data[2,4] = 17.64924845131027
$data[2,5] = 2.07955589531849
$data[2,6] = 2.735997085624268
$data[2,7] = 2.819518945045894
$data[2,8] = 0
$data[2,9] = 0
$data[2,10] = 0
$data[2,11] = 0
$data[2,12] = 0
$data[2,13] = 0
$data[2,14] = 0
$data[2,15] = 14.2082728333152
$data[3,0] = 21.39827057864807
$data[3,1] = 16.86313572588381
$data[3,2] = 4.076836119671445
$data[3,3] = 25.57764346556439
$data[3,4] = 17.49084442949468
$data[3,5] = 2.080842831687721
$data[3,6] = 2.667542724177149
$data[3,7] = 2.769161091993968
$data[3,8] = 0
$data[3,9] = 0
$data[3,10] = 0
$data[3,11] = 0
$data[3,12] = 0
$data[3,13] = 0
$data[3,14] = 0
$data[3,15] = 14.14526399565041
$data[4,0] = 21.32384007689132
$data[4,1] = 16.8353225505169
$data[4,2] = 4.073836145152685
$data[4,3] = 25.49061016046219
$data[4,4] = 17.44754621655228
$data[4,5] = 2.081065407239584
$data[4,6] = 2.655917759783867
$data[4,7] = 2.761577121908168
$data[4,8] = 0
$data[4,9] = 0
$data[4,10] = 0
$data[4,11] = 0
$data[4,12] = 0
$data[4,13] = 0
$data[4,14] = 0
$data[4,15] = 14.12160564649061
$data[5,0] = 21.83321247646425
$data[5,1] = 17.19321626836599
$data[5,2] = 4.117147111988736
$data[5,3] = 26.08648815999262
$data[5,4] = 17.60102890494004
$data[5,5] = 2.079593170369083
$data[5,6] = 2.734610986851889
$data[5,7] = 2.820713870448712
$data[5,8] = 0
$data[5,9] = 0
$data[5,10] = 0
$data[5,11] = 0
$data[5,12] = 0
$data[5,13] = 0
$data[5,14] = 0
$data[5,15] = 14.17126314780701
$data[6,0] = 23.95016255225594
$data[6,1] = 18.65577656276812
$data[6,2] = 4.296175493917771
$data[6,3] = 28.57366895627258
$data[6,4] = 18.33033639464374
$data[6,5] = 2.073275125169522
$data[6,6] = 3.070445310790362
$data[6,7] = 3.073227600650782
$data[6,8] = 0
$data[6,9] = 0
$data[6,10] = 0
$data[6,11] = 0
$data[6,12] = 0
$data[6,13] = 0
$data[6,14] = 0
$data[6,15] = 14.44169348795232
$data[7,0] = 27.63402988064382
$data[7,1] = 21.197725204966
$data[7,2] = 4.614753228655371
$data[7,3] = 32.94140957282742
$data[7,4] = 19.85430217178764
$data[7,5] = 2.06152701904097
$data[7,6] = 3.688970248645137
$data[7,7] = 3.5447492924246
$data[7,8] = 0
$data[7,9] = 0
$data[7,10] = 0
$data[7,11] = 0
$data[7,12] = 0
$data[7,13] = 0
$data[7,14] = 0
$data[7,15] = 15.09616439954949
$data[8,0] = 30.05482127450283
$data[8,1] = 22.928038186333
$data[8,2] = 4.874742813704069
$data[8,3] = 34.9646762972166
$data[8,4] = 20.76645788651782
$data[8,5] = 2.053507625778628
$data[8,6] = 4.0828560900355
$data[8,7] = 3.872324092016694
$data[8,8] = 0
$data[8,9] = 0
$data[8,10] = 0
$data[8,11] = 0
$data[8,12] = 0
$data[8,13] = 0
$data[8,14] = 0
$data[8,15] = 15.44422833104606
$data[9,0] = 31.06280360374938
$data[9,1] = 23.67909720321394
$data[9,2] = 5.319755378876647
$data[9,3] = 28.3645749766292
$data[9,4] = 19.49143900051977
$data[9,5] = 2.051789214662098
$data[9,6] = 4.451631542486619
$data[9,7] = 3.959039045160723
$data[9,8] = 0
$data[9,9] = 0
$data[9,10] = 0
$data[9,11] = 0
$data[9,12] = 0
$data[9,13] = 0
$data[9,14] = 0
$data[9,15] = 14.22126316644512
$data[10,0] = 31.42548808711483
$data[10,1] = 23.9023462235017
$data[10,2] = 5.642827013362487
$data[10,3] = 22.41126104158498
$data[10,4] = 18.26442248682248
$data[10,5] = 2.051750605528799
$data[10,6] = 5.267815255078704
$data[10,7] = 3.968982522168912
$data[10,8] = 0
$data[10,9] = 0
$data[10,10] = 0
$data[10,11] = 0
$data[10,12] = 0
$data[10,13] = 0
$data[10,14] = 0
$data[10,15] = 13.1551535732224
$data[11,0] = 31.32574147658034
$data[11,1] = 23.8023350085204
$data[11,2] = 5.904668838461662
$data[11,3] = 16.48110235038558
$data[11,4] = 16.90881854901085
$data[11,5] = 2.053024761394272
$data[11,6] = 6.295261935674344
$data[11,7] = 3.922986396371375
$data[11,8] = 0
$data[11,9] = 0
$data[11,10] = 0
$data[11,11] = 0
$data[11,12] = 0
$data[11,13] = 0
$data[11,14] = 0
$data[11,15] = 12.08345473812526
$data[12,0] = 31.04587492091814
$data[12,1] = 23.59796749859168
$data[12,2] = 6.064911838226821
$data[12,3] = 12.39948689670872
$data[12,4] = 15.87950193874549
$data[12,5] = 2.054486513402981
$data[12,6] = 7.098170924744764
$data[12,7] = 3.867096236766066
$data[12,8] = 0
$data[12,9] = 0
$data[12,10] = 0
$data[12,11] = 0
$data[12,12] = 0
$data[12,13] = 0
$data[12,14] = 0
$data[12,15] = 11.32485027251851
$data[13,0] = 30.8771560147254
$data[13,1] = 23.49405501822589
$data[13,2] = 6.093011043100423
$data[13,3] = 11.41723285621635
$data[13,4] = 15.5867530759949
$data[13,5] = 2.055159410205087
$data[13,6] = 7.282752320104646
$data[13,7] = 3.841502725540322
$data[13,8] = 0
$data[13,9] = 0
$data[13,10] = 0
$data[13,11] = 0
$data[13,12] = 0
$data[13,13] = 0
$data[13,14] = 0
$data[13,15] = 11.1299896727916
$data[14,0] = 29.90608954675489
$data[14,1] = 22.82702198544147
$data[14,2] = 5.945605588702913
$data[14,3] = 11.28027543478703
$data[14,4] = 15.37602402775759
$data[14,5] = 2.058301267800048
$data[14,6] = 7.001793318481235
$data[14,7] = 3.713436795145757
$data[14,8] = 0
$data[14,9] = 0
$data[14,10] = 0
$data[14,11] = 0
$data[14,12] = 0
$data[14,13] = 0
$data[14,14] = 0
$data[14,15] = 11.15774339611703
$data[15,0] = 29.29981657985098
$data[15,1] = 22.41983437938699
$data[15,2] = 5.742303216938881
$data[15,3] = 13.38635553275821
$data[15,4] = 15.7802112669253
$data[15,5] = 2.059961030497345
$data[15,6] = 6.307931783648366
$data[15,7] = 3.644269112860384
$data[15,8] = 0
$data[15,9] = 0
$data[15,10] = 0
$data[15,11] = 0
$data[15,12] = 0
$data[15,13] = 0
$data[15,14] = 0
$data[15,15] = 11.58299135364639
$data[16,0] = 28.95452436606275
$data[16,1] = 22.16444665626625
$data[16,2] = 5.465780429740172
$data[16,3] = 17.95867208990727
$data[16,4] = 16.78554409291046
$data[16,5] = 2.060395421184307
$data[16,6] = 5.266155076833628
$data[16,7] = 3.619276282305147
$data[16,8] = 0
$data[16,9] = 0
$data[16,10] = 0
$data[16,11] = 0
$data[16,12] = 0
$data[16,13] = 0
$data[16,14] = 0
$data[16,15] = 12.42595252609114
$data[17,0] = 28.85164769548753
$data[17,1] = 22.13903500602907
$data[17,2] = 5.173506130908312
$data[17,3] = 24.24597306498533
$data[17,4] = 18.11461374072026
$data[17,5] = 2.059659617645809
$data[17,6] = 4.278060908810593
$data[17,7] = 3.643304718169249
$data[17,8] = 0
$data[17,9] = 0
$data[17,10] = 0
$data[17,11] = 0
$data[17,12] = 0
$data[17,13] = 0
$data[17,14] = 0
$data[17,15] = 13.49383433671839
$data[18,0] = 29.43636915149335
$data[18,1] = 22.61293314624111
$data[18,2] = 4.829660403423384
$data[18,3] = 34.40385939798826
$data[18,4] = 20.3946872154569
$data[18,5] = 2.0556477065682
$data[18,6] = 3.977166495449919
$data[18,7] = 3.792827574933824
$data[18,8] = 0
$data[18,9] = 0
$data[18,10] = 0
$data[18,11] = 0
$data[18,12] = 0
$data[18,13] = 0
$data[18,14] = 0
$data[18,15] = 15.24019517860266
$data[19,0] = 31.2128873691445
$data[19,1] = 23.87180592760401
$data[19,2] = 4.963924874935277
$data[19,3] = 37.23156527326921
$data[19,4] = 21.41888838623526
$data[19,5] = 2.049227471009077
$data[19,6] = 4.331854164063795
$data[19,7] = 4.05173183932555
$data[19,8] = 0
$data[19,9] = 0
$data[19,10] = 0
$data[19,11] = 0
$data[19,12] = 0
$data[19,13] = 0
$data[19,14] = 0
$data[19,15] = 15.78587279155189
$data[20,0] = 32.32256361032284
$data[20,1] = 24.6138801197208
$data[20,2] = 5.059435747161129
$data[20,3] = 38.57128553915183
$data[20,4] = 22.03675548966449
$data[20,5] = 2.045186385304464
$data[20,6] = 4.540833089149837
$data[20,7] = 4.212966268023601
$data[20,8] = 0
$data[20,9] = 0
$data[20,10] = 0
$data[20,11] = 0
$data[20,12] = 0
$data[20,13] = 0
$data[20,14] = 0
$data[20,15] = 16.11096486298121
$data[21,0] = 31.73534804202202
$data[21,1] = 24.17410756943924
$data[21,2] = 5.001185173529103
$data[21,3] = 37.86180514530616
$data[21,4] = 21.75145750447885
$data[21,5] = 2.047321382313501
$data[21,6] = 4.430018046149989
$data[21,7] = 4.124763455075183
$data[21,8] = 0
$data[21,9] = 0
$data[21,10] = 0
$data[21,11] = 0
$data[21,12] = 0
$data[21,13] = 0
$data[21,14] = 0
$data[21,15] = 15.97437081622241
$data[22,0] = 29.41053052578191
$data[22,1] = 22.52351586716033
$data[22,2] = 4.786763194111646
$data[22,3] = 35.06520157565711
$data[22,4] = 20.60081722598498
$data[22,5] = 2.055542008347234
$data[22,6] = 4.002604939480796
$data[22,7] = 3.79065573902421
$data[22,8] = 0
$data[22,9] = 0
$data[22,10] = 0
$data[22,11] = 0
$data[22,12] = 0
$data[22,13] = 0
$data[22,14] = 0
$data[22,15] = 15.41479316140643
$data[23,0] = 26.68629855093528
$data[23,1] = 20.61659098146711
$data[23,2] = 4.543036311527937
$data[23,3] = 31.81296454054986
$data[23,4] = 19.36590674883738
$data[23,5] = 2.064667127147425
$data[23,6] = 3.52510738843625
$data[23,7] = 3.422478011629868
$data[23,8] = 0
$data[23,9] = 0
$data[23,10] = 0
$data[23,11] = 0
$data[23,12] = 0
$data[23,13] = 0
$data[23,14] = 0
$data[23,15] = 14.85062811068665

$ws.Range("B2:Q25").Value = $data
